$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 99
$ws.Range("D2").Value = 72

$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 103.5

$ws.Range("C4").Value = 96
$ws.Range("D4").Value = 62

$ws.Range("C5").Value = 127
$ws.Range("D5").Value = 104.5

$ws.Range("C6").Value = 116
$ws.Range("D6").Value = 102.5

$ws.Range("C7").Value = 124
$ws.Range("D7").Value = 99

$ws.Range("C8").Value = 63
$ws.Range("D8").Value = 39.5

$ws.Range("C9").Value = 74
$ws.Range("D9").Value = 69

$ws.Range("C10").Value = 91
$ws.Range("D10").Value = 91

$ws.Range("C11").Value = 64
$ws.Range("D11").Value = 39.5

$ws.Range("C12").Value = 88.8
